$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New attendance rows to append (same ID/Name/Date, new check-in/out times)
$rows = @(
    @("1446896", "Asif Newaz", "2025-01-25", "20:45:12"),
    @("1446896", "Asif Newaz", "2025-01-25", "23:02:05"),
    @("1446896", "Asif Newaz", "2025-01-25", "23:03:56")
)

$startRow = 3
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $rowRange = $ws.Range("A$r`:D$r")

    # Force these to be stored as text (like the existing rows), so values
    # such as "1446896" or "2025-01-25" aren't reinterpreted as a number/date.
    $rowRange.NumberFormat = "@"

    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("D$r").Value = $data[3]

    # Drop the formatting we applied so the new cells end up styled the
    # same as the rest of the sheet (no explicit style/number format).
    $rowRange.ClearFormats()
}
